# Adding noise filter to matching algorithm
# Insert a new parameter row ("CorrelationMinimum") above the
# "IsotopicPercentage" row, shifting the rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5 - everything currently at/after row 5 shifts down.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row with the new parameter's data.
$ws.Range("A5").Value = "CorrelationMinimum"
$ws.Range("B5").Value = 0.7
$ws.Range("C5").Value = "Everytime"
$ws.Range("D5").Value = "The minimum correlation value to consider when generating the trelliscope display"

# Update the active selection to match the edited workbook.
$ws.Range("B6").Select()
